# Add order quantities to the PolyTag BOM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / fill in a couple of DIGIKEY part numbers ---
$ws.Range("G10").Value = "490-1318-1-ND"
$ws.Range("G14").Value = "160-1446-1-ND"

# --- New header cells (match the existing bold/bordered header look) ---
$ws.Range("K1").Value = "Needed"
$ws.Range("L1").Value = "Order 2015-07-30"
$ws.Range("M1").Value = "Extra"
$ws.Range("J1").Copy()
$ws.Range("K1:M1").PasteSpecial(-4122)

# --- Order-quantity helper cell (multiplier used by the "Needed" formula) ---
$ws.Range("H32").Value = "Wanted"
$ws.Range("I32").Value = 36

# --- Per-row Needed / Order / Extra values ---
# Needed = Qty * $I$32 ; Extra = Order - Needed
$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("K$r").Formula = "=A$r*`$I`$32"
    $ws.Range("M$r").Formula = "=L$r-K$r"
}

$ws.Range("L2").Value = 150
$ws.Range("L5").Value = 900
$ws.Range("L8").Value = 83
$ws.Range("L9").Value = 100
$ws.Range("L10").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("L12").Value = 83
$ws.Range("L13").Value = 83
$ws.Range("L14").Value = 500
$ws.Range("L15").Value = 40
$ws.Range("L16").Value = 100
$ws.Range("L17").Value = 83
$ws.Range("L18").Value = 76
$ws.Range("L19").Value = 153
$ws.Range("L20").Value = 500
$ws.Range("L21").Value = 750
$ws.Range("L22").Value = 80
$ws.Range("L23").Value = 500
$ws.Range("L24").Value = 500
$ws.Range("L26").Value = 36
$ws.Range("L27").Value = 100
$ws.Range("L28").Value = 100
$ws.Range("L29").Value = 50

# --- Conditional formatting: flag shortfalls (Extra < 0) in red ---
$rng = $ws.Range("M2:M29")
$fc = $rng.FormatConditions.Add(1, 6, "0")
$fc.Font.Color = 255

# --- Column widths: column F shrinks back to the G width, K/L get a sensible width ---
$ws.Columns.Item(6).ColumnWidth = 30.98
$ws.Columns.Item(11).ColumnWidth = 15.46
$ws.Columns.Item(12).ColumnWidth = 15.46

# --- Leave the cursor where the last edit happened ---
$ws.Range("L15").Select()
